$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.533.40'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.76%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.540.93'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.37%  '

$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '312.13'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.40%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '100.39'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.45%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.568'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.07%  '

$ws.Range("E8").Value = '  +0.07%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.528'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.90%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.84'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.85%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0802'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.15%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.38'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.42%  '

$ws.Range("E13").Value = '  -0.44%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.933.61'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.70%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.88'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +5.81%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.553.66'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.05%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.831'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.93%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.569.30'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.68%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.77'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.41%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0952'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.17%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.28'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.92%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.96'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.73%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '243.64'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.78%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.92'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.12%  '

$ws.Range("E25").Value = '  +0.31%  '

$ws.Range("B26").Value = 'Dai'
$ws.Range("C26").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.03%  '

$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '26.28'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.41%  '

$ws.Range("B28").Value = 'InjectiveProtocol'
$ws.Range("C28").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '40.34'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.68%  '

$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.33'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.27%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '10.08'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.01%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '158.12'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.21%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.67'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.41%  '

$ws.Range("E33").Value = '  +14.70%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0801'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.58%  '

$ws.Range("B35").Value = 'WEMIXToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.62'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.13%  '

$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.05'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.12%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.19'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.08%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.17'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.53%  '

$ws.Range("E39").Value = '  -1.60%  '

$ws.Range("E40").Value = '  -0.51%  '

$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.17'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +9.82%  '

$ws.Range("B42").Value = 'EnergySwap'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '22.04'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.95%  '

$ws.Range("E43").Value = '  +0.28%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.31'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.30%  '

$ws.Range("E45").Value = '  -1.78%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.965.47'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.55%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.85'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.08%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.789.06'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.67%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '80.82'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.28%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.192'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.43%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '72.88'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.09%  '
